# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 16:39"

# Row 4 (e.g. India)
$ws.Range("B4").Value = 6488716
$ws.Range("C4").Value = 3141
$ws.Range("D4").Value = 3759134
$ws.Range("E4").Value = 2535986
$ws.Range("G4").Value = 62
$ws.Range("H4").Value = 193596

# Row 5
$ws.Range("B5").Value = 4313129
$ws.Range("C5").Value = 35545
$ws.Range("D5").Value = 3352316
$ws.Range("E5").Value = 887708
$ws.Range("G5").Value = 289
$ws.Range("H5").Value = 73105

# Row 13
$ws.Range("D13").Value = 366590
$ws.Range("E13").Value = 111238
$ws.Range("G13").Value = 50
$ws.Range("H13").Value = 10179

# Row 24
$ws.Range("B24").Value = 254168
$ws.Range("C24").Value = 543
$ws.Range("E24").Value = 17761
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 9407

# Row 47
$ws.Range("B47").Value = 73208
$ws.Range("C47").Value = 177
$ws.Range("E47").Value = 604
$ws.Range("G47").Value = 5
$ws.Range("H47").Value = 721

# Row 51
$ws.Range("B51").Value = 60895
$ws.Range("C51").Value = 388
$ws.Range("D51").Value = 43146
$ws.Range("E51").Value = 15903
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 1846

# Row 62
$ws.Range("E62").Value = 5120
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 2017

# Row 64
$ws.Range("B64").Value = 44281
$ws.Range("C64").Value = 388
$ws.Range("E64").Value = 2329
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 358

# Row 67
$ws.Range("B67").Value = 37557
$ws.Range("C67").Value = 139
$ws.Range("D67").Value = 34965
$ws.Range("E67").Value = 2040
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 552

# Row 69
$ws.Range("B69").Value = 31994
$ws.Range("C69").Value = 53
$ws.Range("E69").Value = 630
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 727

# Row 73
$ws.Range("D73").Value = 17779
$ws.Range("E73").Value = 9394

# Row 92
$ws.Range("B92").Value = 11560
$ws.Range("C92").Value = 39
$ws.Range("E92").Value = 1948

# Row 99
$ws.Range("B99").Value = 8860
$ws.Range("C99").Value = 36
$ws.Range("D99").Value = 7650
$ws.Range("E99").Value = 1140

# Row 124
$ws.Range("B124").Value = 3900
$ws.Range("C124").Value = 124
$ws.Range("D124").Value = 1817
$ws.Range("E124").Value = 2037
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 46

# Row 165
$ws.Range("B165").Value = 1054
$ws.Range("C165").Value = 5
$ws.Range("D165").Value = 868
$ws.Range("E165").Value = 151

# Row 175
$ws.Range("B175").Value = 503
$ws.Range("C175").Value = 6
$ws.Range("E175").Value = 266
